# LV Activities - 29 July
# Update the primary-attendee / capital-provider test company name on the
# "Company" sheet, and leave that sheet active/selected at A2 (mirrors a
# tester opening the Company tab, clicking A2 and retyping its value).

$wb = $excel.ActiveWorkbook

$companySheet = $wb.Worksheets.Item("Company")
$companySheet.Activate()

$companySheet.Range("A2").Value = "Capital Provider Test Company"
$companySheet.Range("A2").Select()
